$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at N:P (existing N,O,P shift right to Q,R,S)
$ws.Range("N1:P1").EntireColumn.Insert()

# New column headers
$ws.Range("N1").Value = "stat_u"
$ws.Range("O1").Value = "syst_u"
$ws.Range("P1").Value = "syst_c"

# Row 2 (non-shared formulas, matching the single-cell formulas in the diff)
$ws.Range("N2").Formula = "=K2*J2/100"
$ws.Range("O2").Formula = "=L2*J2/100"
$ws.Range("P2").Formula = "=M2*J2/100"

# Rows 3-9 (fill down as in the diff's shared formula ranges)
$ws.Range("N3:N9").Formula = "=K3*J3/100"
$ws.Range("O3:O9").Formula = "=L3*J3/100"
$ws.Range("P3:P9").Formula = "=M3*J3/100"

# Update selection to match the edited workbook's saved cursor position
$ws.Activate()
$ws.Range("P2:P9").Select()

$wb.Save()
